$p = $ppt.ActivePresentation

# --- Slide 2: reposition the "Picture 4" photo slightly, and remove the
#     arrow connector + its "15 min walk" caption that pointed at it.
$s2 = $p.Slides.Item(2)

$picture4 = $s2.Shapes.Item("Picture 4")
$picture4.Left = 653.82665
$picture4.Top  = 325.3756

$s2.Shapes.Item("TextBox 9").Delete()
$s2.Shapes.Item("Straight Arrow Connector 8").Delete()

# --- Remove the duplicate title/intro slide (slide 12).
$p.Slides.Item(12).Delete()
